# Add a "Save" column (H) to the s_vals sheet, mirroring the style of the
# existing "sum" header (G1) and filling in the per-row save flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 ("sum") onto H1, then set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new column's data values for the existing rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
